$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 10 (shift existing data rows 10-34 down to 13-37)
$ws.Rows("10:12").Insert()

# Common (constant) values for this subset, same across all rows
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103003
$categoria = "Damasco"

# New week's data: 2021-12-16 (serial 44546), Castle Brite
$rows = @(
    @{ Row=10; Fecha=44546; Variedad="Castle Brite"; Calidad="Especial"; Volumen=240; PMin=21500; PMax=22000; PProm=21750; Unidad="`$/caja 18 kilos"; Origen="Región Metropolitana"; PrecioKg=1208; KgUnidad=18 },
    @{ Row=11; Fecha=44546; Variedad="Castle Brite"; Calidad="Primera";  Volumen=300; PMin=19500; PMax=20000; PProm=19750; Unidad="`$/caja 18 kilos"; Origen="Región Metropolitana"; PrecioKg=1097; KgUnidad=18 },
    @{ Row=12; Fecha=44546; Variedad="Castle Brite"; Calidad="Segunda";  Volumen=300; PMin=15500; PMax=16000; PProm=15750; Unidad="`$/caja 18 kilos"; Origen="Región Metropolitana"; PrecioKg=875;  KgUnidad=18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}

Write-Host "done"
